# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# fdc22035-... source file has moved from "In Translation" to
# "Ready for handoff" (a new handoff was generated for it), and
# refreshes the associated handoff timestamps on every sheet.

$wb = $excel.ActiveWorkbook

# ---- "zh-cn" sheet ---------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-29 12:14:41"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ---- "de-de" sheet -----------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-29 12:14:45"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25

# ---- "Overview" sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-29 12:14:45"
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
